# Append: 2025-11-12 06:35 JST
#
# A new scrape run found one additional listing that now ranks above the
# existing "priority score 10" rows. It gets inserted immediately before
# the current row 18 on the "ランサーズ" sheet, pushing the old rows
# 18-20 down to 19-21. Every row's "取得日時" (retrieved-at) timestamp is
# refreshed to the new scrape time, including the rows that shifted down
# and the brand-new row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2025-11-12 06:35:53"

# --- Shift the bottom three data rows down by one (18->19, 19->20, 20->21) ---
# Copy plain values only, bottom-up so nothing is clobbered before it is
# read. The cells already carry the "Hyperlink" style/format from the
# template, so it rides along with the value automatically; only the
# brand-new last row needs a fresh hyperlink registered below.
for ($r = 20; $r -ge 18; $r--) {
    for ($c = 1; $c -le 7; $c++) {
        $val = $ws.Cells.Item($r, $c).Value()
        $ws.Cells.Item($r + 1, $c).Value = $val
    }
}

# --- Write the brand-new listing into row 18 ---
$ws.Cells.Item(18, 1).Value = $newTimestamp
$ws.Cells.Item(18, 2).Value = "【急募】Wartalesの武器アイコンとモデルを日本刀に差し替え"
$ws.Cells.Item(18, 3).Value = "システム開発"
$ws.Cells.Item(18, 4).Value = "10,000 円 ~ 20,000 円 / 固定"
$ws.Cells.Item(18, 5).Value = "期限情報なし"
$ws.Cells.Item(18, 6).Value = "https://www.lancers.jp/work/detail/5432425"
$ws.Cells.Item(18, 7).Value = 10

# --- Register a hyperlink for the now-last row (21), which previously
#     didn't exist. Existing hyperlinks on rows 18-20 already moved down
#     together with their cells, so they don't need touching. ---
$ws.Hyperlinks.Add($ws.Range("F21"), "https://www.lancers.jp/work/detail/5432055", [Type]::Missing, [Type]::Missing, [Type]::Missing)
$ws.Range("F21").Style = "Hyperlink"

# --- Refresh the "取得日時" timestamp on every data row ---
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 1).Value = $newTimestamp
}
